$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "323.22"
Set-TextValue $ws.Range("E2") "-2.48%"

# Row 3
Set-TextValue $ws.Range("D3") "39.63"
Set-TextValue $ws.Range("E3") "-1.73%"

# Row 4
Set-TextValue $ws.Range("D4") "5.878"
Set-TextValue $ws.Range("E4") "11.38%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08031"
Set-TextValue $ws.Range("E5") "-0.71%"

# Row 6
Set-TextValue $ws.Range("D6") "4.574"
Set-TextValue $ws.Range("E6") "0.82%"

# Row 7
Set-TextValue $ws.Range("D7") "8.670"
Set-TextValue $ws.Range("E7") "0.20%"

# Row 8
Set-TextValue $ws.Range("D8") "1.932"
Set-TextValue $ws.Range("E8") "0.24%"

# Row 9
Set-TextValue $ws.Range("D9") "2.948"
Set-TextValue $ws.Range("E9") "-0.32%"

# Row 10
Set-TextValue $ws.Range("D10") "0.9297"
Set-TextValue $ws.Range("E10") "-0.77%"

# Row 11
Set-TextValue $ws.Range("D11") "0.1262"
Set-TextValue $ws.Range("E11") "-6.01%"

# Row 12
Set-TextValue $ws.Range("D12") "0.1971"
Set-TextValue $ws.Range("E12") "0.38%"

# Row 13
Set-TextValue $ws.Range("D13") "8.716"
Set-TextValue $ws.Range("E13") "34.19%"

# Row 14
Set-TextValue $ws.Range("D14") "0.09130"
Set-TextValue $ws.Range("E14") "-0.02%"

# Row 15
Set-TextValue $ws.Range("D15") "0.03573"
Set-TextValue $ws.Range("E15") "1.98%"

# Row 16
Set-TextValue $ws.Range("E16") "9.27%"

# Row 17
Set-TextValue $ws.Range("D17") "0.001299"
Set-TextValue $ws.Range("E17") "-5.66%"

# Row 18
Set-TextValue $ws.Range("D18") "0.006318"
Set-TextValue $ws.Range("E18") "-1.57%"

# Row 19
Set-TextValue $ws.Range("E19") "-0.34%"

# Row 20
Set-TextValue $ws.Range("D20") "0.3536"
Set-TextValue $ws.Range("E20") "0.39%"

# Row 21
Set-TextValue $ws.Range("D21") "0.1371"
Set-TextValue $ws.Range("E21") "3.34%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2450"
Set-TextValue $ws.Range("E22") "-4.99%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04412"
Set-TextValue $ws.Range("E23") "-0.81%"

# Row 24
Set-TextValue $ws.Range("E24") "3.16%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004376"
Set-TextValue $ws.Range("E25") "2.19%"

# Row 26
Set-TextValue $ws.Range("E26") "-11.71%"

# Row 39
Set-TextValue $ws.Range("D39") "0.02506"
Set-TextValue $ws.Range("E39") "-0.02%"

# Row 40
Set-TextValue $ws.Range("D40") "0.05281"
Set-TextValue $ws.Range("E40") "1.62%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007429"
Set-TextValue $ws.Range("E41") "-3.05%"

# Row 42
Set-TextValue $ws.Range("D42") "0.009611"
Set-TextValue $ws.Range("E42") "4.58%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1405"
Set-TextValue $ws.Range("E43") "-1.37%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002117"
Set-TextValue $ws.Range("E44") "0.05%"

# Row 45
Set-TextValue $ws.Range("D45") "0.009987"
Set-TextValue $ws.Range("E45") "21.51%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006728"
Set-TextValue $ws.Range("E46") "1.72%"

# Row 47
Set-TextValue $ws.Range("E47") "-0.11%"

# Row 48
Set-TextValue $ws.Range("D48") "0.003006"
Set-TextValue $ws.Range("E48") "-10.13%"

# Row 49
Set-TextValue $ws.Range("D49") "0.002291"
Set-TextValue $ws.Range("E49") "-7.73%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "-0.11%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "-0.11%"
